$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report date header (D1)
$ws.Range("D1").Value = "24/03/2023"

# Update forecast/activity values in column D
$ws.Range("D3").Value = 287.3
$ws.Range("D4").Value = 314.6
$ws.Range("D5").Value = 334.1
$ws.Range("D6").Value = 393.9
$ws.Range("D7").Value = 265.2
$ws.Range("D8").Value = 252.2
$ws.Range("D9").Value = 380.9
$ws.Range("D10").Value = 310.7
$ws.Range("D11").Value = 257.4
$ws.Range("D12").Value = 291.2
$ws.Range("D13").Value = 234
$ws.Range("D14").Value = 148.2
$ws.Range("D18").Value = 8
$ws.Range("D19").Value = 35
$ws.Range("D20").Value = 36
$ws.Range("D21").Value = 28
$ws.Range("D22").Value = 33
$ws.Range("D23").Value = 18
$ws.Range("D24").Value = 24
$ws.Range("D25").Value = 20
$ws.Range("D26").Value = 32
$ws.Range("D27").Value = 21
$ws.Range("D28").Value = 26
$ws.Range("D29").Value = 13
$ws.Range("D30").Value = 9
$ws.Range("D34").Value = 32
$ws.Range("D35").Value = 68
$ws.Range("D36").Value = 148
$ws.Range("D37").Value = 221
$ws.Range("D38").Value = 189
$ws.Range("D39").Value = 180
$ws.Range("D40").Value = 153
$ws.Range("D41").Value = 137
$ws.Range("D42").Value = 140
$ws.Range("D43").Value = 134
$ws.Range("D44").Value = 145
$ws.Range("D45").Value = 146
$ws.Range("D46").Value = 157
$ws.Range("D47").Value = 122
$ws.Range("D48").Value = 80
$ws.Range("D49").Value = 51
$ws.Range("D50").Value = 19
$ws.Range("D51").Value = 12
$ws.Range("D52").Value = 79
$ws.Range("D53").Value = 138
$ws.Range("D54").Value = 149
$ws.Range("D55").Value = 150
$ws.Range("D56").Value = 141
$ws.Range("D57").Value = 101
$ws.Range("D58").Value = 105
$ws.Range("D59").Value = 121
$ws.Range("D60").Value = 99
$ws.Range("D61").Value = 126
$ws.Range("D62").Value = 95
$ws.Range("D63").Value = 66
$ws.Range("D68").Value = 26
$ws.Range("D69").Value = 45
$ws.Range("D70").Value = 44
$ws.Range("D71").Value = 45
$ws.Range("D72").Value = 34
$ws.Range("D73").Value = 38
$ws.Range("D74").Value = 41
$ws.Range("D75").Value = 40
$ws.Range("D76").Value = 38
$ws.Range("D77").Value = 41
$ws.Range("D78").Value = 42
$ws.Range("D79").Value = 21
$ws.Range("D80").Value = 14
$ws.Range("D81").Value = 8
$ws.Range("D82").Value = 3
$ws.Range("D83").Value = 2
$ws.Range("D84").Value = 12
$ws.Range("D85").Value = 2
$ws.Range("D86").Value = 3
$ws.Range("D87").Value = 2
$ws.Range("D88").Value = 4
$ws.Range("D89").Value = 6
$ws.Range("D91").Value = 107
$ws.Range("D92").Value = 72
$ws.Range("D93").Value = 165
$ws.Range("D94").Value = 155
$ws.Range("D95").Value = 151
$ws.Range("D96").Value = 147
$ws.Range("D97").Value = 131
$ws.Range("D98").Value = 160
$ws.Range("D99").Value = 135
$ws.Range("D100").Value = 135
$ws.Range("D101").Value = 148
$ws.Range("D102").Value = 156
$ws.Range("D103").Value = 109
$ws.Range("D104").Value = 78
$ws.Range("D105").Value = 48
$ws.Range("D106").Value = 32
$ws.Range("D107").Value = 12
$ws.Range("D109").Value = 22
$ws.Range("D110").Value = 25
$ws.Range("D111").Value = 26
$ws.Range("D112").Value = 20
$ws.Range("D113").Value = 14
$ws.Range("D114").Value = 22
$ws.Range("D115").Value = 15
$ws.Range("D116").Value = 20
$ws.Range("D117").Value = 17
$ws.Range("D118").Value = 16
$ws.Range("D119").Value = 12
$ws.Range("D120").Value = 5
$ws.Range("D126").Value = 0
$ws.Range("D129").Value = 2
$ws.Range("D131").Value = 0
$ws.Range("D133").Value = 1
$ws.Range("D134").Value = 0
$ws.Range("D135").Value = 0
$ws.Range("D140").Value = 6
$ws.Range("D141").Value = 15
$ws.Range("D142").Value = 13
$ws.Range("D143").Value = 13
$ws.Range("D144").Value = 14
$ws.Range("D146").Value = 13
$ws.Range("D147").Value = 11
$ws.Range("D148").Value = 15
$ws.Range("D149").Value = 17
$ws.Range("D150").Value = 19
$ws.Range("D151").Value = 12
$ws.Range("D152").Value = 6
$ws.Range("D153").Value = 4
$ws.Range("D154").Value = 4
